# Vaccination diagram: make the diagram "a little less thick".
#
# This slimmer look is achieved by:
#   1. Reducing the bend point of the red elbow/bent connector that runs
#      from "Susceptible" down to "Recovered" (shape "Connector: Elbow 2"),
#      pulling it further away from the bottom edge.
#   2. Moving the "Vaccination" textbox up to follow the connector's new
#      bend point.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- 1) Shrink the adjustment of the bent connector ("Connector: Elbow 2") ---
# OOXML: <a:gd name="adj1" fmla="val 5280795"/> -> <a:gd name="adj1" fmla="val 2931323"/>
# The PowerPoint object model expresses this as a percentage-like value
# (raw OOXML fmla value / 100000).
$connector = $s.Shapes.Item("Connector: Elbow 2")
$connector.Adjustments.Item(1) = 29.31323

# --- 2) Move the "Vaccination" textbox up to match ---
# OOXML: <a:off x="3692069" y="2578388"/> -> <a:off x="3692069" y="2365948"/>
# Shape.Top/Left are expressed in points (EMU / 12700).
$vaccinationBox = $s.Shapes.Item("TextBox 11")
$vaccinationBox.Top = 2365948 / 12700
